$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 brand-new rows before the current row 891, shifting the existing
# rows 891:980 down to 894:983 (dimension grows from T980 to T983).
$ws.Rows("891:893").Insert()

# ---- New row 891 ----
$ws.Cells.Item(891, 1).Value = 10
$ws.Cells.Item(891, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(891, 3).Value = "La Araucanía"
$ws.Cells.Item(891, 4).Value = 45132
$ws.Cells.Item(891, 5).Value = 9
$ws.Cells.Item(891, 6).Value = "Fruta"
$ws.Cells.Item(891, 7).Value = 100102
$ws.Cells.Item(891, 8).Value = "Cítricos"
$ws.Cells.Item(891, 9).Value = 100102004
$ws.Cells.Item(891, 10).Value = "Mandarina"
$ws.Cells.Item(891, 11).Value = "Clementina"
$ws.Cells.Item(891, 12).Value = "Especial"
$ws.Cells.Item(891, 13).Value = 300
$ws.Cells.Item(891, 14).Value = 10000
$ws.Cells.Item(891, 15).Value = 10000
$ws.Cells.Item(891, 16).Value = 10000
$ws.Cells.Item(891, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(891, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(891, 19).Value = 1000
$ws.Cells.Item(891, 20).Value = 10

# ---- New row 892 ----
$ws.Cells.Item(892, 1).Value = 10
$ws.Cells.Item(892, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(892, 3).Value = "La Araucanía"
$ws.Cells.Item(892, 4).Value = 45132
$ws.Cells.Item(892, 5).Value = 9
$ws.Cells.Item(892, 6).Value = "Fruta"
$ws.Cells.Item(892, 7).Value = 100102
$ws.Cells.Item(892, 8).Value = "Cítricos"
$ws.Cells.Item(892, 9).Value = 100102004
$ws.Cells.Item(892, 10).Value = "Mandarina"
$ws.Cells.Item(892, 11).Value = "Murcott"
$ws.Cells.Item(892, 12).Value = "Primera"
$ws.Cells.Item(892, 13).Value = 250
$ws.Cells.Item(892, 14).Value = 16000
$ws.Cells.Item(892, 15).Value = 16000
$ws.Cells.Item(892, 16).Value = 16000
$ws.Cells.Item(892, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(892, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(892, 19).Value = 889
$ws.Cells.Item(892, 20).Value = 18

# ---- New row 893 ----
$ws.Cells.Item(893, 1).Value = 10
$ws.Cells.Item(893, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(893, 3).Value = "La Araucanía"
$ws.Cells.Item(893, 4).Value = 45132
$ws.Cells.Item(893, 5).Value = 9
$ws.Cells.Item(893, 6).Value = "Fruta"
$ws.Cells.Item(893, 7).Value = 100102
$ws.Cells.Item(893, 8).Value = "Cítricos"
$ws.Cells.Item(893, 9).Value = 100102004
$ws.Cells.Item(893, 10).Value = "Mandarina"
$ws.Cells.Item(893, 11).Value = "Tangerina Kara"
$ws.Cells.Item(893, 12).Value = "Primera"
$ws.Cells.Item(893, 13).Value = 80
$ws.Cells.Item(893, 14).Value = 17000
$ws.Cells.Item(893, 15).Value = 17000
$ws.Cells.Item(893, 16).Value = 17000
$ws.Cells.Item(893, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(893, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(893, 19).Value = 944
$ws.Cells.Item(893, 20).Value = 18

# Apply the date-number-format (numFmtId 165, "YYYY-MM-DD HH:MM:SS") used by
# every other row's "Fecha" column to the D cells of the three new rows —
# matching the format already carried onto D891:D893 by the row Insert().
$ws.Range("D891:D893").NumberFormat = $ws.Range("D894").NumberFormat
